$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12 (old rows 12-16 shift down to 13-17),
# preserving the existing column styles (A: s=1, B: s=2).
$ws.Rows("12:12").Insert()

# New event e011 "Deployment" content, inserted right after e010 (row 11).
$eventId = "e011"

$body = @'
<Bold>e011 Deployment</Bold> 
<InlineUIContainer><Button Content='r4.41' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine your tank&apos;s deployment from the 
<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table:  
<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@

$ws.Range("A12").Value = $eventId
$ws.Range("B12").Value = $body

# Match the row height used for similar-length entries (e.g. row 9).
$ws.Rows("12:12").RowHeight = 99.85

# Restore the active selection shown in the saved workbook.
$ws.Range("B11").Select()
